# Rename the worksheet and update the active selection to match
# what was recorded after comparing the FEC / FX option to xVA Lite.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet1 -> constant_vol_surface
$ws.Name = "constant_vol_surface"

# Move the active cell selection from D8 to E14
$ws.Range("E14").Select()
